$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D (Price) updates - NumberFormat forced to Text ("@") where the new
# value would otherwise be auto-parsed as a number, so the cell keeps storing
# the exact original text (e.g. trailing zeros like "299.73", "0.520").
$ws.Range("D2").Value = '42.960.38'
$ws.Range("D3").Value = '2.293.44'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.73'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.71'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.520'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.93'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.58'
$ws.Range("D15").Value = '2.650.46'
$ws.Range("D16").Value = '2.297.76'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.784'
$ws.Range("D18").Value = '42.884.54'
$ws.Range("D20").Value = '0.0₃0911'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.79'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.61'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.80'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.08'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.05'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.99'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.73'
$ws.Range("D42").Value = '2.000.83'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.20'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.32'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.81'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.95'
$ws.Range("D49").Value = '2.518.01'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '73.27'

# Column E (Volume/1h change %) updates - plain text, padded with two spaces
# on each side to match the existing formatting convention in the sheet.
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("E6").Value = '  -2.77%  '
$ws.Range("E7").Value = '  +1.97%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("E11").Value = '  -0.75%  '
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("E14").Value = '  -1.48%  '
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("E19").Value = '  -0.02%  '
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  -3.35%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -1.13%  '
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("E30").Value = '  -1.16%  '
$ws.Range("E31").Value = '  -4.80%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +2.76%  '
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("E37").Value = '  +0.20%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("E39").Value = '  -1.49%  '
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("E41").Value = '  -1.99%  '
$ws.Range("E42").Value = '  +0.95%  '
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("E44").Value = '  -1.88%  '
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("E48").Value = '  -2.48%  '
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  +3.33%  '
$ws.Range("E51").Value = '  -2.06%  '
